$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale per-row notes in column E (rows 8 and 10): these were
# "server down 03.07.2019" and "too slow" comments that are no longer needed.
$ws.Range("E8").ClearContents()
$ws.Range("E10").ClearContents()

# Make room for six new Skosmos instances as rows 11-16, inheriting the
# un-bold "data row" formatting used by the rows above (style matches C2/D2).
$ws.Rows("11:16").Insert()

# Add the new name/url pairs (context/timeout counters left at 0), entering
# values in the same order as the authored edit (url before name for some
# rows) so shared strings line up the same way.
$ws.Range("A11").Value = "GACS"
$ws.Range("B11").Value = "http://artemide.art.uniroma2.it/skosmos/"
$ws.Range("B12").Value = "http://51.15.194.251/Skosmos/"
$ws.Range("B13").Value = "https://voc.uni-ak.ac.at/skosmos"
$ws.Range("A12").Value = "51.15.194.251"
$ws.Range("A13").Value = "UAAV"
$ws.Range("B14").Value = "http://skosmos.linkeddata.ch/"
$ws.Range("A14").Value = "HTW Chur"
$ws.Range("B15").Value = "https://isl.ics.forth.gr/apollonis-federated-thesaurus/"
$ws.Range("A15").Value = "FORTH"
$ws.Range("B16").Value = "https://vocabulaires.irstea.fr/skosmos/"
$ws.Range("A16").Value = "Irstea"

$ws.Range("C11:D16").Value = 0

# Match the final selection / active cell recorded in the worksheet view.
[void]$ws.Range("E14").Select()
